# Auto-generated cell updates applying the diff to the cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the cells we are about to touch to remain plain text so that
# values such as '1.010' or '0.5280' are not reinterpreted as numbers.
$cellsToUpdate = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'E26', 'D27', 'D28', 'E28', 'D29', 'E29', 'D30', 'E30', 'E31', 'D32', 'E32', 'E33', 'D34', 'E34', 'D35', 'E35', 'D36', 'E36', 'D37', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'D46', 'E46', 'B47', 'C47', 'D47', 'E47', 'B48', 'C48', 'D48', 'E48', 'B49', 'C49', 'D49', 'E49', 'B50', 'C50', 'D50', 'E50', 'D51', 'E51')
foreach ($addr in $cellsToUpdate) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '26.308.44'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.667.50'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '219.89'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').Value = '0.5280'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.2647'
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').Value = '0.06364'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = '20.92'
$ws.Range('E10').Value = '  +1.84%  '
$ws.Range('D11').Value = '0.07839'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '4.522'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').Value = '1.669.78'
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').Value = '1.897.08'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').Value = '0.5604'
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').Value = '0.0₅8100'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '65.76'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').Value = '26.334.09'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = '1.010'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').Value = '4.720'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('D21').Value = '199.68'
$ws.Range('E21').Value = '  +3.62%  '
$ws.Range('D22').Value = '10.27'
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').Value = '6.051'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').Value = '1.011'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = '146.51'
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('E26').Value = '  -0.99%  '
$ws.Range('D27').Value = '7.243'
$ws.Range('D28').Value = '16.16'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = '1.510'
$ws.Range('E29').Value = '  +2.65%  '
$ws.Range('D30').Value = '0.05884'
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').Value = '3.513'
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('D34').Value = '1.598'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').Value = '0.9635'
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('D36').Value = '2.821'
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('D37').Value = '2.433'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = '0.5794'
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('D39').Value = '0.01615'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '5.951'
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('D41').Value = '1.076.37'
$ws.Range('E41').Value = '  +2.68%  '
$ws.Range('D42').Value = '0.8586'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('D43').Value = '1.010'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').Value = '102.94'
$ws.Range('E44').Value = '  -1.54%  '
$ws.Range('D45').Value = '1.807.74'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('D46').Value = '58.43'
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').Value = '1.012'
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.4418'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '8.085'
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₈103'
$ws.Range('E50').Value = '  -3.83%  '
$ws.Range('D51').Value = '0.05149'
$ws.Range('E51').Value = '  -0.33%  '
